$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "保險" (insurance) -> sheet6.xml
# Columns: A index, B company, C name, D owner, E category, F property_category,
#          G date, H legislator_name, I legislator_id, J source_file, K index
# ---------------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item(6)

$wsIns.Cells.Item(1, 2).Value = "company"
$wsIns.Cells.Item(1, 3).Value = "name"
$wsIns.Cells.Item(1, 4).Value = "owner"
$wsIns.Cells.Item(1, 5).Value = "property_category"
$wsIns.Cells.Item(1, 6).Value = "category"
$wsIns.Cells.Item(1, 7).Value = "date"
$wsIns.Cells.Item(1, 8).Value = "legislator_name"
$wsIns.Cells.Item(1, 9).Value = "legislator_id"
$wsIns.Cells.Item(1, 10).Value = "source_file"
$wsIns.Cells.Item(1, 11).Value = "index"

$insRows = @(
    @(84, "富邦人壽", "20LPL安泰分红終身壽險", "蘇震清"),
    @(85, "富邦人壽", "20LPL安泰分紅終身壽險", "廖靖汝"),
    @(86, "富邦人壽", "20LPL安泰分紅終身壽險", "蘇震清"),
    @(87, "富邦人壽", "20LPL安泰分红終身壽險", "蘇震清"),
    @(88, "富邦人壽", "安泰喬壽還本終身壽險", "蘇震清"),
    @(89, "富邦人壽", "安泰喬壽還本終身壽險", "蘇震清")
)

for ($insIdx = 0; $insIdx -lt $insRows.Count; $insIdx++) {
    $insR = $insIdx + 2
    $insRow = $insRows[$insIdx]
    $wsIns.Cells.Item($insR, 2).Value = $insRow[1]
    $wsIns.Cells.Item($insR, 3).Value = $insRow[2]
    $wsIns.Cells.Item($insR, 4).Value = $insRow[3]
    $wsIns.Cells.Item($insR, 5).Value = "insurance"
    $wsIns.Cells.Item($insR, 6).Value = "normal"
    $wsIns.Cells.Item($insR, 7).Value = "2012-04-30"
    $wsIns.Cells.Item($insR, 8).Value = "蘇震清"
    $wsIns.Cells.Item($insR, 9).Value = 1718
    $wsIns.Cells.Item($insR, 10).Value = "tmp16a71"
    $wsIns.Cells.Item($insR, 11).Value = $insRow[0]
}

# ---------------------------------------------------------------------------
# Sheet "債務" (debt) -> sheet7.xml
# Columns: A index, B species, C debtor, D owner, E total, F register_date,
#          G register_reason, H property_category, I category, J date,
#          K legislator_name, L legislator_id, M source_file, N index
# ---------------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item(7)

$wsDebt.Cells.Item(1, 2).Value = "species"
$wsDebt.Cells.Item(1, 3).Value = "debtor"
$wsDebt.Cells.Item(1, 4).Value = "owner"
$wsDebt.Cells.Item(1, 5).Value = "total"
$wsDebt.Cells.Item(1, 6).Value = "register_date"
$wsDebt.Cells.Item(1, 7).Value = "register_reason"
$wsDebt.Cells.Item(1, 8).Value = "property_category"
$wsDebt.Cells.Item(1, 9).Value = "category"
$wsDebt.Cells.Item(1, 10).Value = "date"
$wsDebt.Cells.Item(1, 11).Value = "legislator_name"
$wsDebt.Cells.Item(1, 12).Value = "legislator_id"
$wsDebt.Cells.Item(1, 13).Value = "source_file"
$wsDebt.Cells.Item(1, 14).Value = "index"

$wsDebt.Cells.Item(2, 2).Value = "綜合存款存摺(透支）"
$wsDebt.Cells.Item(2, 3).Value = "蘇震清"
$wsDebt.Cells.Item(2, 4).Value = "台灣銀行屏東分行屏東縣屏東市中山路"
$wsDebt.Cells.Item(2, 5).Value = 4766183
$wsDebt.Cells.Item(2, 6).Value = "87年11月03日"
$wsDebt.Cells.Item(2, 7).Value = "貸款中期擔保放款"
$wsDebt.Cells.Item(2, 8).Value = "debt"
$wsDebt.Cells.Item(2, 9).Value = "normal"
$wsDebt.Cells.Item(2, 10).Value = "2012-04-30"
$wsDebt.Cells.Item(2, 11).Value = "蘇震清"
$wsDebt.Cells.Item(2, 12).Value = 1718
$wsDebt.Cells.Item(2, 13).Value = "tmp16a71"
$wsDebt.Cells.Item(2, 14).Value = 99

$wsDebt.Cells.Item(3, 2).Value = "長期擔保放款"
$wsDebt.Cells.Item(3, 3).Value = "廖靖汝"
$wsDebt.Cells.Item(3, 4).Value = "台灣土地銀行屏東縣屛東市逢甲路"
$wsDebt.Cells.Item(3, 5).Value = 616200
$wsDebt.Cells.Item(3, 6).Value = "87年11月03日"
$wsDebt.Cells.Item(3, 7).Value = "房貸"
$wsDebt.Cells.Item(3, 8).Value = "debt"
$wsDebt.Cells.Item(3, 9).Value = "normal"
$wsDebt.Cells.Item(3, 10).Value = "2012-04-30"
$wsDebt.Cells.Item(3, 11).Value = "蘇震清"
$wsDebt.Cells.Item(3, 12).Value = 1718
$wsDebt.Cells.Item(3, 13).Value = "tmp16a71"
$wsDebt.Cells.Item(3, 14).Value = 100

$wsDebt.Cells.Item(4, 2).Value = "綜合存款存摺(透支）"
$wsDebt.Cells.Item(4, 3).Value = "蘇震清"
$wsDebt.Cells.Item(4, 4).Value = "台灣銀行屏東分行屏東縣屏東市中山路"
$wsDebt.Cells.Item(4, 5).Value = 27993
$wsDebt.Cells.Item(4, 6).Value = "100年03月28日"
$wsDebt.Cells.Item(4, 7).Value = "貸款治家成長貸款"
$wsDebt.Cells.Item(4, 8).Value = "debt"
$wsDebt.Cells.Item(4, 9).Value = "normal"
$wsDebt.Cells.Item(4, 10).Value = "2012-04-30"
$wsDebt.Cells.Item(4, 11).Value = "蘇震清"
$wsDebt.Cells.Item(4, 12).Value = 1718
$wsDebt.Cells.Item(4, 13).Value = "tmp16a71"
$wsDebt.Cells.Item(4, 14).Value = 101

# ---------------------------------------------------------------------------
# Sheet "事業投資" (business investment) -> sheet8.xml
# Columns: A index, B owner, C company, D address, E total, F register_date,
#          G register_reason, H property_category, I category, J date,
#          K legislator_name, L legislator_id, M source_file, N index
# ---------------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item(8)

$wsInv.Cells.Item(1, 2).Value = "owner"
$wsInv.Cells.Item(1, 3).Value = "company"
$wsInv.Cells.Item(1, 4).Value = "address"
$wsInv.Cells.Item(1, 5).Value = "total"
$wsInv.Cells.Item(1, 6).Value = "register_date"
$wsInv.Cells.Item(1, 7).Value = "register_reason"
$wsInv.Cells.Item(1, 8).Value = "property_category"
$wsInv.Cells.Item(1, 9).Value = "category"
$wsInv.Cells.Item(1, 10).Value = "date"
$wsInv.Cells.Item(1, 11).Value = "legislator_name"
$wsInv.Cells.Item(1, 12).Value = "legislator_id"
$wsInv.Cells.Item(1, 13).Value = "source_file"
$wsInv.Cells.Item(1, 14).Value = "index"

$wsInv.Cells.Item(2, 1).Value = 106
$wsInv.Cells.Item(2, 2).Value = "廖靖汝"
$wsInv.Cells.Item(2, 3).Value = "南島休閒育樂股份有限公司"
$wsInv.Cells.Item(2, 4).Value = "高雄市精富路148號"
$wsInv.Cells.Item(2, 5).Value = 1000000
$wsInv.Cells.Item(2, 6).Value = "95年08月23日"
$wsInv.Cells.Item(2, 7).Value = "投資"
$wsInv.Cells.Item(2, 8).Value = "investment"
$wsInv.Cells.Item(2, 9).Value = "normal"
$wsInv.Cells.Item(2, 10).Value = "2012-04-30"
$wsInv.Cells.Item(2, 11).Value = "蘇震清"
$wsInv.Cells.Item(2, 12).Value = 1718
$wsInv.Cells.Item(2, 13).Value = "tmp16a71"
$wsInv.Cells.Item(2, 14).Value = 106
